$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "5.22", "0.540") are stored verbatim as text instead
# of being coerced into floating point numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '57.563.11'
$ws.Range("E2").Value = '  +5.93%  '

# Row 3
$ws.Range("D3").Value = '2.337.04'
$ws.Range("E3").Value = '  +2.89%  '

# Row 4
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.30%  '

# Row 5
$ws.Range("D5").Value = '522.06'
$ws.Range("E5").Value = '  +4.58%  '

# Row 6
$ws.Range("D6").Value = '135.22'
$ws.Range("E6").Value = '  +4.56%  '

# Row 7
$ws.Range("D7").Value = '0.994'
$ws.Range("E7").Value = '  -0.35%  '

# Row 8
$ws.Range("D8").Value = '0.540'
$ws.Range("E8").Value = '  +2.63%  '

# Row 9
$ws.Range("D9").Value = '2.367.37'
$ws.Range("E9").Value = '  +3.80%  '

# Row 10
$ws.Range("E10").Value = '  +9.14%  '

# Row 11
$ws.Range("E11").Value = '  +0.98%  '

# Row 12
$ws.Range("D12").Value = '5.22'
$ws.Range("E12").Value = '  +5.81%  '

# Row 13
$ws.Range("E13").Value = '  +2.74%  '

# Row 14
$ws.Range("D14").Value = '24.06'
$ws.Range("E14").Value = '  +3.90%  '

# Row 15
$ws.Range("D15").Value = '2.755.96'
$ws.Range("E15").Value = '  +3.11%  '

# Row 16
$ws.Range("D16").Value = '57.267.72'
$ws.Range("E16").Value = '  +5.42%  '

# Row 17
$ws.Range("D17").Value = '0.0000136'
$ws.Range("E17").Value = '  +5.14%  '

# Row 18
$ws.Range("D18").Value = '2.363.19'
$ws.Range("E18").Value = '  +4.17%  '

# Row 19
$ws.Range("D19").Value = '10.63'
$ws.Range("E19").Value = '  +3.46%  '

# Row 20
$ws.Range("D20").Value = '4.30'
$ws.Range("E20").Value = '  +3.59%  '

# Row 21
$ws.Range("D21").Value = '323.17'
$ws.Range("E21").Value = '  +6.14%  '

# Row 22
$ws.Range("E22").Value = '  +5.62%  '

# Row 23
$ws.Range("D23").Value = '0.998'

# Row 24
$ws.Range("D24").Value = '61.45'
$ws.Range("E24").Value = '  +1.30%  '

# Row 25
$ws.Range("D25").Value = '0.993'
$ws.Range("E25").Value = '  -0.40%  '

# Row 26
$ws.Range("E26").Value = '  +7.34%  '

# Row 27
$ws.Range("D27").Value = '7.81'
$ws.Range("E27").Value = '  +5.87%  '

# Row 28
$ws.Range("D28").Value = '172.21'
$ws.Range("E28").Value = '  -1.76%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0750'
$ws.Range("E29").Value = '  +6.38%  '

# Row 30
$ws.Range("E30").Value = '  +12.27%  '

# Row 31
$ws.Range("D31").Value = '6.35'
$ws.Range("E31").Value = '  +5.62%  '

# Row 32
$ws.Range("D32").Value = '1.69'
$ws.Range("E32").Value = '  +4.75%  '

# Row 33
$ws.Range("E33").Value = '  +3.58%  '

# Row 34
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.04%  '

# Row 35
$ws.Range("D35").Value = '0.960'
$ws.Range("E35").Value = '  +0.69%  '

# Row 36
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -0.01%  '

# Row 37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '1.27'
$ws.Range("E37").Value = '  +5.47%  '

# Row 38
$ws.Range("D38").Value = '4.05'
$ws.Range("E38").Value = '  +8.41%  '

# Row 39
$ws.Range("E39").Value = '  +8.88%  '

# Row 40
$ws.Range("D40").Value = '37.65'
$ws.Range("E40").Value = '  +4.26%  '

# Row 41
$ws.Range("D41").Value = '0.384'
$ws.Range("E41").Value = '  +2.16%  '

# Row 42
$ws.Range("D42").Value = '140.15'
$ws.Range("E42").Value = '  +11.81%  '

# Row 43
$ws.Range("D43").Value = '3.62'
$ws.Range("E43").Value = '  +6.95%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '5.27'
$ws.Range("E44").Value = '  +3.67%  '

# Row 45
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = '278.71'
$ws.Range("E45").Value = '  +13.18%  '

# Row 46
$ws.Range("E46").Value = '  +4.06%  '

# Row 47
$ws.Range("D47").Value = '0.0933'

# Row 48
$ws.Range("E48").Value = '  +3.99%  '

# Row 49
$ws.Range("D49").Value = '0.0217'
$ws.Range("E49").Value = '  +5.91%  '

# Row 50
$ws.Range("D50").Value = '0.384'
$ws.Range("E50").Value = '  +2.25%  '

# Row 51
$ws.Range("D51").Value = '17.08'
$ws.Range("E51").Value = '  +5.15%  '

# Remove the temporary text-format override so the cell style stays
# identical to the original workbook (no explicit style index).
$ws.Range("D2:E51").ClearFormats()
